# Generate Report for Handback
#
# The localization status report is regenerated after the handback for
# file "533b14e5-b481-45f0-9339-8a423ddbac09.md" completes successfully
# (it is now "in sync" with en-US instead of merely "Ready for handoff").
#
# This updates:
#   - Overview sheet: per-language status cells for that file's row
#   - zh-cn sheet: Status / Latest Handback DateTime / Error Detail for that row
#   - de-de sheet: Status / Latest Handback DateTime / Error Detail for that row

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"

# ---- Overview sheet ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $statusHandedBack
$wsOverview.Range("F3").Value = $statusHandedBack

# ---- zh-cn sheet ----
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $statusHandedBack
$wsZhCn.Range("K3").Value = "2016-08-31 13:01:37"
$wsZhCn.Range("P3").Value = ""
$wsZhCn.Columns.Item(16).ColumnWidth = 12.9

# ---- de-de sheet ----
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $statusHandedBack
$wsDeDe.Range("K3").Value = "2016-08-31 13:01:45"
$wsDeDe.Range("P3").Value = ""
$wsDeDe.Columns.Item(16).ColumnWidth = 12.9
